$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Update input values on Sheet1 (drives the recalculated formula results) ---
$ws1.Range("B3").Value = 18.742215440623902
$ws1.Range("B6").Value = 1.2949999999999999
$ws1.Range("B7").Value = 3.1192901530031798
$ws1.Range("B8").Value = 0.19400000000000001
$ws1.Range("B12").Value = 73.456879999999998
$ws1.Range("B13").Value = 2

# --- Add a new reference/annotation cell (leading apostrophe forces text + quotePrefix, like typing '-this...' into Excel) ---
$ws1.Range("F12").Value = "'-this is a matt reference"

# --- View state: zoom in and move the selection ---
$ws1.Select()
$excel.ActiveWindow.Zoom = 220
$ws1.Range("A17").Select()

$excel.CalculateFull()
